$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "4:1 Conditioned" (C) and "4:1 Unconditioned" (D) observation
# counts for every data row.
for ($r = 2; $r -le 78; $r++) {
    $cVal = $ws.Cells.Item($r, 3).Value2
    $dVal = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 3).Value = $dVal
    $ws.Cells.Item($r, 4).Value = $cVal
}

# Update the sheet view: scroll back to the top, zoom to 114%, and select F7.
$ws.Range("F7").Select() | Out-Null
$excel.ActiveWindow.Zoom = 114
